$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 104, pushing the existing rows 104-114 down to 105-115.
$ws.Rows("104:104").Insert()

# Populate the newly inserted row 104 with the new weekly price record.
$ws.Range("A104").Value = 11
$ws.Range("B104").Value = "Vega Monumental Concepción"
$ws.Range("C104").Value = "Bíobío"
$ws.Range("D104").Value = 45077
$ws.Range("E104").Value = 8
$ws.Range("F104").Value = 100112037
$ws.Range("G104").Value = "Cebollín"
$ws.Range("H104").Value = "Sin especificar"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 40
$ws.Range("K104").Value = 3500
$ws.Range("L104").Value = 4000
$ws.Range("M104").Value = 3750
$ws.Range("N104").Value = "`$/paquete 36 unidades"
$ws.Range("O104").Value = "Región Metropolitana"
$ws.Range("P104").Value = 104
$ws.Range("Q104").Value = 36
$ws.Range("R104").Value = "Hortaliza"
